$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '48.058.74'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.499.27'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.70'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.82'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.523'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.71'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.03'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0804'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.09'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.889.01'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.496.92'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.833'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.880.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.04'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.95'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.65'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0934'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.22'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '272.88'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +9.61%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.74'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.75%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.141'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.85'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.26'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.15'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.29'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0775'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.94'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.57'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.66%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '121.56'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.111'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.94%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.28'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.21'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.004.55'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.16'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.88'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.18'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.86'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.08%  '
